# Update attendee/view counts (column F) on the "展览" and "全部类型" sheets
# to reflect the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 7366
$wsExhibit.Range("F4").Value = 5771
$wsExhibit.Range("F5").Value = 87
$wsExhibit.Range("F13").Value = 81
$wsExhibit.Range("F14").Value = 656
$wsExhibit.Range("F15").Value = 446
$wsExhibit.Range("F20").Value = 65

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 7366
$wsAll.Range("F4").Value = 5771
$wsAll.Range("F5").Value = 87
$wsAll.Range("F13").Value = 81
$wsAll.Range("F15").Value = 446
$wsAll.Range("F20").Value = 65
